# Update database: drop oldest quarter (column D, فصل دوم منتهی به 1399/06)
# and append newest quarter (فصل چهارم منتهی به 1401/12) as the new last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the oldest quarter column (D). This shifts every later quarter
#    one column to the left (E->D, F->E, ... M->L) and drops the obsolete
#    "فصل دوم منتهی به 1399/06" period entirely.
$ws.Range("D1").EntireColumn.Delete()

# 2. The newest quarter column is now M (previously unused / off the
#    right edge). Populate its header (period) and publish-date labels,
#    and give it the same (wider) column width used for the other
#    year-end ("فصل چهارم") quarter columns, e.g. column E.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30 (2)"
$ws.Range("M1").ColumnWidth = $ws.Range("E1").ColumnWidth

# 3. The publish date recorded for the quarter that is now in column I
#    (فصل چهارم منتهی به 1400/12) was corrected/updated.
$ws.Range("I9").Value = "1402-02-30 (8)"

# 4. Fill in the financial figures for the new quarter (column M).
$ws.Range("M11").Value = 6440315
$ws.Range("M12").Value = -3184789
$ws.Range("M13").Value = 3255526
$ws.Range("M14").Value = -192499
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = -3578
$ws.Range("M17").Value = 3059449
$ws.Range("M18").Value = -617000
$ws.Range("M19").Value = 77132
$ws.Range("M20").Value = 2519581
$ws.Range("M21").Value = 4962
$ws.Range("M22").Value = 2524543
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 2524543
$ws.Range("M25").Value = 1122
$ws.Range("M26").Value = 2250000
$ws.Range("M27").Value = 1122
